# Updates PETR4 "Robo PUT" option-chain sheet with refreshed market data
# (TIR %, Strike VS Cot. %, Prob. Exec., Negocios/volume, and a couple of
# Robo PUT status flips) - "novo bot no telegram" data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = -0.63
$ws.Range("G5").Value = -0.62
$ws.Range("G6").Value = -0.6
$ws.Range("I9").Value = 200
$ws.Range("E12").Value = 0.03
$ws.Range("F12").Value = 0.0022
$ws.Range("G12").Value = -0.49
$ws.Range("I12").Value = 10000
$ws.Range("G13").Value = -0.47
$ws.Range("H13").Value = 0
$ws.Range("G14").Value = -0.45
$ws.Range("H14").Value = 0.01
$ws.Range("G15").Value = -0.43
$ws.Range("G16").Value = -0.42
$ws.Range("H16").Value = 0.04
$ws.Range("G17").Value = -0.41
$ws.Range("H17").Value = 0.05
$ws.Range("G18").Value = -0.4
$ws.Range("H18").Value = 0.08
$ws.Range("G19").Value = -0.39
$ws.Range("H19").Value = 0.11
$ws.Range("I19").Value = 7000
$ws.Range("G20").Value = -0.38
$ws.Range("H20").Value = 0.15
$ws.Range("G21").Value = -0.37
$ws.Range("H21").Value = 0.2
$ws.Range("G22").Value = -0.36
$ws.Range("H22").Value = 0.27
$ws.Range("G23").Value = -0.35
$ws.Range("H23").Value = 0.36
$ws.Range("G24").Value = -0.34
$ws.Range("H24").Value = 0.48
$ws.Range("G25").Value = -0.33
$ws.Range("H25").Value = 0.62
$ws.Range("H26").Value = 1.01
$ws.Range("H27").Value = 1.27
$ws.Range("H28").Value = 1.59
$ws.Range("G29").Value = -0.29
$ws.Range("H29").Value = 1.96
$ws.Range("E30").Value = 0.13
$ws.Range("F30").Value = 0.0068
$ws.Range("G30").Value = -0.28
$ws.Range("H30").Value = 2.39
$ws.Range("I30").Value = 4100
$ws.Range("E31").Value = 0.13
$ws.Range("F31").Value = 0.0068
$ws.Range("G31").Value = -0.27
$ws.Range("H31").Value = 2.9
$ws.Range("I31").Value = 34800
$ws.Range("G32").Value = -0.26
$ws.Range("H32").Value = 3.48
$ws.Range("E33").Value = 0.15
$ws.Range("F33").Value = 0.0076
$ws.Range("G33").Value = -0.25
$ws.Range("H33").Value = 4.15
$ws.Range("I33").Value = 3000
$ws.Range("G34").Value = -0.24
$ws.Range("H34").Value = 4.9
$ws.Range("E35").Value = 0.18
$ws.Range("F35").Value = 0.0089
$ws.Range("G35").Value = -0.23
$ws.Range("H35").Value = 5.76
$ws.Range("I35").Value = 100
$ws.Range("E36").Value = 0.19
$ws.Range("F36").Value = 0.009299999999999999
$ws.Range("G36").Value = -0.22
$ws.Range("H36").Value = 6.71
$ws.Range("I36").Value = 500
$ws.Range("G37").Value = -0.21
$ws.Range("H37").Value = 7.76
$ws.Range("E38").Value = 0.22
$ws.Range("F38").Value = 0.0105
$ws.Range("G38").Value = -0.2
$ws.Range("H38").Value = 8.93
$ws.Range("I38").Value = 4200
$ws.Range("G39").Value = -0.19
$ws.Range("H39").Value = 10.2
$ws.Range("I39").Value = 100
$ws.Range("E40").Value = 0.27
$ws.Range("F40").Value = 0.0126
$ws.Range("G40").Value = -0.18
$ws.Range("H40").Value = 11.58
$ws.Range("I40").Value = 100
$ws.Range("G41").Value = -0.17
$ws.Range("H41").Value = 13.07
$ws.Range("E42").Value = 0.32
$ws.Range("F42").Value = 0.0145
$ws.Range("G42").Value = -0.16
$ws.Range("H42").Value = 14.67
$ws.Range("I42").Value = 68600
$ws.Range("E43").Value = 0.34
$ws.Range("F43").Value = 0.0153
$ws.Range("G43").Value = -0.15
$ws.Range("H43").Value = 16.37
$ws.Range("I43").Value = 100
$ws.Range("E44").Value = 0.38
$ws.Range("F44").Value = 0.0169
$ws.Range("G44").Value = -0.14
$ws.Range("H44").Value = 18.18
$ws.Range("I44").Value = 3500
$ws.Range("G45").Value = -0.13
$ws.Range("H45").Value = 20.09
$ws.Range("D46").Value = "Aguardar"
$ws.Range("E46").Value = 0.46
$ws.Range("F46").Value = 0.02
$ws.Range("G46").Value = -0.12
$ws.Range("H46").Value = 22.09
$ws.Range("I46").Value = 60300
$ws.Range("E47").Value = 0.5
$ws.Range("F47").Value = 0.0215
$ws.Range("G47").Value = -0.11
$ws.Range("H47").Value = 24.17
$ws.Range("I47").Value = 1100
$ws.Range("D48").Value = "Montar"
$ws.Range("E48").Value = 0.57
$ws.Range("F48").Value = 0.0243
$ws.Range("G48").Value = -0.1
$ws.Range("H48").Value = 26.34
$ws.Range("I48").Value = 1500
$ws.Range("E49").Value = 0.61
$ws.Range("F49").Value = 0.0257
$ws.Range("G49").Value = -0.09
$ws.Range("H49").Value = 28.57
$ws.Range("I49").Value = 30100
$ws.Range("E50").Value = 0.6899999999999999
$ws.Range("F50").Value = 0.0288
$ws.Range("G50").Value = -0.09
$ws.Range("H50").Value = 30.86
$ws.Range("I50").Value = 61300
$ws.Range("E51").Value = 0.74
$ws.Range("F51").Value = 0.0305
$ws.Range("G51").Value = -0.08
$ws.Range("H51").Value = 33.21
$ws.Range("I51").Value = 400
$ws.Range("E52").Value = 0.9
$ws.Range("F52").Value = 0.0364
$ws.Range("G52").Value = -0.06
$ws.Range("H52").Value = 38.03
$ws.Range("I52").Value = 11600
$ws.Range("E53").Value = 0.99
$ws.Range("F53").Value = 0.0396
$ws.Range("G53").Value = -0.05
$ws.Range("H53").Value = 40.47
$ws.Range("I53").Value = 9000
$ws.Range("E54").Value = 1.06
$ws.Range("F54").Value = 0.042
$ws.Range("G54").Value = -0.04
$ws.Range("H54").Value = 42.93
$ws.Range("I54").Value = 100
$ws.Range("E55").Value = 1.19
$ws.Range("F55").Value = 0.0467
$ws.Range("G55").Value = -0.03
$ws.Range("H55").Value = 45.39
$ws.Range("I55").Value = 3600
$ws.Range("G56").Value = -0.03
$ws.Range("H56").Value = 45.39
$ws.Range("E57").Value = 1.3
$ws.Range("F57").Value = 0.0505
$ws.Range("G57").Value = -0.02
$ws.Range("H57").Value = 47.85
$ws.Range("I57").Value = 1100
$ws.Range("E58").Value = 1.38
$ws.Range("F58").Value = 0.0531
$ws.Range("G58").Value = -0.01
$ws.Range("H58").Value = 50.29
$ws.Range("I58").Value = 42100
$ws.Range("E59").Value = 1.53
$ws.Range("F59").Value = 0.0583
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = 52.7
$ws.Range("I59").Value = 2500
$ws.Range("E60").Value = 1.65
$ws.Range("F60").Value = 0.0623
$ws.Range("G60").Value = 0.01
$ws.Range("H60").Value = 55.09
$ws.Range("I60").Value = 8700
$ws.Range("E61").Value = 1.79
$ws.Range("F61").Value = 0.0669
$ws.Range("G61").Value = 0.02
$ws.Range("H61").Value = 57.43
$ws.Range("I61").Value = 13000
$ws.Range("E62").Value = 1.91
$ws.Range("F62").Value = 0.0707
$ws.Range("G62").Value = 0.03
$ws.Range("H62").Value = 59.72
$ws.Range("I62").Value = 6200
$ws.Range("E63").Value = 2.25
$ws.Range("F63").Value = 0.0818
$ws.Range("G63").Value = 0.05
$ws.Range("H63").Value = 64.15000000000001
$ws.Range("I63").Value = 500
$ws.Range("G64").Value = 0.05
$ws.Range("H64").Value = 64.15000000000001
$ws.Range("G65").Value = 0.06
$ws.Range("H65").Value = 66.27
$ws.Range("E66").Value = 2.53
$ws.Range("F66").Value = 0.09039999999999999
$ws.Range("G66").Value = 0.07000000000000001
$ws.Range("H66").Value = 68.31999999999999
$ws.Range("I66").Value = 64200
$ws.Range("G67").Value = 0.09
$ws.Range("H67").Value = 72.19
$ws.Range("G68").Value = 0.09
$ws.Range("H68").Value = 72.19
$ws.Range("G69").Value = 0.1
$ws.Range("H69").Value = 74.02
$ws.Range("G70").Value = 0.11
$ws.Range("H70").Value = 75.77
$ws.Range("G71").Value = 0.12
$ws.Range("H71").Value = 79.02
$ws.Range("G72").Value = 0.14
$ws.Range("H72").Value = 81.95
$ws.Range("G73").Value = 0.16
$ws.Range("H73").Value = 84.56999999999999
$ws.Range("G74").Value = 0.16
$ws.Range("H74").Value = 84.56999999999999
$ws.Range("G75").Value = 0.17
$ws.Range("H75").Value = 85.77
$ws.Range("G76").Value = 0.18
$ws.Range("H76").Value = 86.89
$ws.Range("G77").Value = 0.2
$ws.Range("H77").Value = 88.92
$ws.Range("G78").Value = 0.21
$ws.Range("H78").Value = 89.84
$ws.Range("G79").Value = 0.22
$ws.Range("H79").Value = 90.69
$ws.Range("G80").Value = 0.24
$ws.Range("H80").Value = 92.22
$ws.Range("G81").Value = 0.26
$ws.Range("H81").Value = 93.54000000000001
